$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 17 with hero default / "at large" entry
# Set D17 (#CCCCCC) before B17 (在野) so the shared-string table gets the
# new unique strings appended in the same order as the reference edit
# (index 39 = "#CCCCCC", index 40 = "在野").
$ws.Range("A17").Value = 99
$ws.Range("C17").Value = 100020
$ws.Range("D17").Value = "#CCCCCC"
$ws.Range("B17").Value = "在野"

# Copy style formatting from row 16 to keep consistent look
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)

# Update the selected cell as recorded in the workbook view
$ws.Range("D12").Select()
